$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = "0606_091803.png"
$ws.Range("I2").Value = "06-06 09:18:03 setText execution was Passed"

$ws.Range("G3").Value = "0606_091804.png"
$ws.Range("I3").Value = "06-06 09:18:04 click execution was Passed"
